$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Western Sahara" (ESH) row. It sat at row 182 (Vietnam row 181 +1,
# Yemen previously at 183). Deleting the whole row shifts every following row
# up by one and keeps all the A:B pairings intact.
$ws.Rows.Item(182).Delete()

# Keep the hidden AutoFilter-database defined name in sync with the new,
# one-row-shorter data range.
$names = $wb.Names
$names.Item("_xlnm._FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$B`$184"

# Restore the view: scroll so row 168 is at the top and select the new last
# data row, A185 (the first empty row right below the shrunk table).
$win = $excel.ActiveWindow
$win.ScrollRow = 168
$win.ScrollColumn = 1
$ws.Range("A185").Select() | Out-Null
